$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update class label from "12 RPL A" to "12 A" for students 1 and 4 (rows 2 and 5)
$ws.Range("O2").Value = "12 A"
$ws.Range("O5").Value = "12 A"

# Clear the class value for student 2 (row 3), fixing export error when siswa kelas is empty
$ws.Range("O3").ClearContents()
